# Insert a new weekly data row at row 38 (shifting existing rows 38-84 down to 39-85)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 5
$ws.Range("B38").Value = "Macroferia Regional de Talca"
$ws.Range("C38").Value = "Maule"
$ws.Range("D38").Value = 45100
$ws.Range("E38").Value = 7
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100107
$ws.Range("H38").Value = "Otros"
$ws.Range("I38").Value = 100107001
$ws.Range("J38").Value = "Caqui"
$ws.Range("K38").Value = "Mankaki"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 150
$ws.Range("N38").Value = 18000
$ws.Range("O38").Value = 18000
$ws.Range("P38").Value = 18000
$ws.Range("Q38").Value = "`$/caja 18 kilos granel"
$ws.Range("R38").Value = "Provincia de Curicó"
$ws.Range("S38").Value = 1000
$ws.Range("T38").Value = 18
